$wb = $excel.ActiveWorkbook

# --- Sheet: Summary ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B2").Value = 0.7294007490636704
$wsSummary.Range("C2").Value = 0.9359430604982206
$wsSummary.Range("D2").Value = 0.4925093632958801
$wsSummary.Range("E2").Value = 0.645398773006135
$wsSummary.Range("F2").Value = 0.5440628878775341
$wsSummary.Range("G2").Value = 0.5016506492553737
$wsSummary.Range("H2").Value = 0.7294007490636705
$wsSummary.Range("I2").Value = 263
$wsSummary.Range("J2").Value = 18
$wsSummary.Range("K2").Value = 516
$wsSummary.Range("L2").Value = 271

# --- Sheet: Classification Report ---
$wsClass = $wb.Worksheets.Item("Classification Report")
$wsClass.Range("B2").Value = 0.6556543837357052
$wsClass.Range("C2").Value = 0.9662921348314607
$wsClass.Range("D2").Value = 0.781226343679031

$wsClass.Range("B3").Value = 0.9359430604982206
$wsClass.Range("C3").Value = 0.4925093632958801
$wsClass.Range("D3").Value = 0.645398773006135

$wsClass.Range("B4").Value = 0.7294007490636704
$wsClass.Range("C4").Value = 0.7294007490636704
$wsClass.Range("D4").Value = 0.7294007490636704
$wsClass.Range("E4").Value = 0.7294007490636704

$wsClass.Range("B5").Value = 0.7957987221169629
$wsClass.Range("C5").Value = 0.7294007490636705
$wsClass.Range("D5").Value = 0.713312558342583

$wsClass.Range("B6").Value = 0.7957987221169629
$wsClass.Range("C6").Value = 0.7294007490636704
$wsClass.Range("D6").Value = 0.7133125583425831

# --- Sheet: Confusion Matrix ---
$wsConf = $wb.Worksheets.Item("Confusion Matrix")
$wsConf.Range("B2").Value = 516
$wsConf.Range("C2").Value = 18
$wsConf.Range("B3").Value = 271
$wsConf.Range("C3").Value = 263
